$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force percentage-looking text cells to remain plain text (avoid Excel auto-converting
# "NN%" into a numeric percentage value), matching the source data which stores these
# as literal strings.
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H41").NumberFormat = "@"

$ws.Range("E2").Value = "2026-02-21 20:18:33"
$ws.Range("E3").Value = "2026-02-21 20:18:36"
$ws.Range("O3").Value = "2.0 °C"
$ws.Range("E4").Value = "2026-02-21 20:18:38"
$ws.Range("O4").Value = "9.6 °C"
$ws.Range("E5").Value = "2026-02-21 20:18:41"
$ws.Range("E6").Value = "2026-02-21 20:18:44"
$ws.Range("E7").Value = "2026-02-21 20:18:46"
$ws.Range("E8").Value = "2026-02-21 20:18:49"
$ws.Range("E9").Value = "2026-02-21 20:18:52"
$ws.Range("H9").Value = "54%"
$ws.Range("N9").Value = "7.6 °C 19:59 TU"
$ws.Range("O9").Value = "13.6 °C"
$ws.Range("E10").Value = "2026-02-21 20:18:54"
$ws.Range("O10").Value = "8.6 °C"
$ws.Range("E11").Value = "2026-02-21 20:18:57"
$ws.Range("H11").Value = "51%"
$ws.Range("O11").Value = "9.1 °C"
$ws.Range("E12").Value = "2026-02-21 20:19:00"
$ws.Range("O12").Value = "12.9 °C"
$ws.Range("E13").Value = "2026-02-21 20:19:02"
$ws.Range("H13").Value = "61%"
$ws.Range("K13").Value = "16.0 MJ/m2"
$ws.Range("E14").Value = "2026-02-21 20:19:05"
$ws.Range("H14").Value = "69%"
$ws.Range("O14").Value = "11.5 °C"
$ws.Range("E15").Value = "2026-02-21 20:19:07"
$ws.Range("H15").Value = "54%"
$ws.Range("N15").Value = "6.6 °C 19:39 TU"
$ws.Range("O15").Value = "13.4 °C"
$ws.Range("E16").Value = "2026-02-21 20:19:10"
$ws.Range("E17").Value = "2026-02-21 20:19:12"
$ws.Range("E18").Value = "2026-02-21 20:19:15"
$ws.Range("O18").Value = "8.8 °C"
$ws.Range("E19").Value = "2026-02-21 20:19:18"
$ws.Range("E20").Value = "2026-02-21 20:19:20"
$ws.Range("E21").Value = "2026-02-21 20:19:23"
$ws.Range("E22").Value = "2026-02-21 20:19:26"
$ws.Range("E23").Value = "2026-02-21 20:19:28"
$ws.Range("K23").Value = "16.0 MJ/m2"
$ws.Range("E24").Value = "2026-02-21 20:19:31"
$ws.Range("E25").Value = "2026-02-21 20:19:33"
$ws.Range("E26").Value = "2026-02-21 20:19:36"
$ws.Range("J26").Value = "1027.3 hPa"
$ws.Range("O26").Value = "9.7 °C"
$ws.Range("E27").Value = "2026-02-21 20:19:39"
$ws.Range("E28").Value = "2026-02-21 20:19:41"
$ws.Range("E29").Value = "2026-02-21 20:19:44"
$ws.Range("E30").Value = "2026-02-21 20:19:47"
$ws.Range("H30").Value = "66%"
$ws.Range("J30").Value = "1029.2 hPa"
$ws.Range("O30").Value = "11.7 °C"
$ws.Range("E31").Value = "2026-02-21 20:19:49"
$ws.Range("O31").Value = "12.2 °C"
$ws.Range("E32").Value = "2026-02-21 20:19:52"
$ws.Range("H32").Value = "79%"
$ws.Range("O32").Value = "5.4 °C"
$ws.Range("E33").Value = "2026-02-21 20:19:55"
$ws.Range("J33").Value = "1030.3 hPa"
$ws.Range("E34").Value = "2026-02-21 20:19:57"
$ws.Range("N34").Value = "-0.2 °C 19:38 TU"
$ws.Range("O34").Value = "4.6 °C"
$ws.Range("E35").Value = "2026-02-21 20:20:00"
$ws.Range("H35").Value = "55%"
$ws.Range("J35").Value = "1030.9 hPa"
$ws.Range("K35").Value = "16.4 MJ/m2"
$ws.Range("O35").Value = "7.7 °C"
$ws.Range("E36").Value = "2026-02-21 20:20:02"
$ws.Range("H36").Value = "57%"
$ws.Range("O36").Value = "13.4 °C"
$ws.Range("E37").Value = "2026-02-21 20:20:05"
$ws.Range("J37").Value = "1031.4 hPa"
$ws.Range("O37").Value = "5.9 °C"
$ws.Range("E38").Value = "2026-02-21 20:20:08"
$ws.Range("H38").Value = "72%"
$ws.Range("O38").Value = "9.8 °C"
$ws.Range("E39").Value = "2026-02-21 20:20:11"
$ws.Range("H39").Value = "34%"
$ws.Range("I39").Value = "1.1 mm"
$ws.Range("K39").Value = "16.2 MJ/m2"
$ws.Range("M39").Value = "5.1 °C 13:33 TU"
$ws.Range("O39").Value = "2.2 °C"
$ws.Range("E40").Value = "2026-02-21 20:20:13"
$ws.Range("J40").Value = "1030.5 hPa"
$ws.Range("O40").Value = "8.9 °C"
$ws.Range("E41").Value = "2026-02-21 20:20:16"
$ws.Range("H41").Value = "67%"
$ws.Range("K41").Value = "15.4 MJ/m2"
$ws.Range("E42").Value = "2026-02-21 20:20:18"
$ws.Range("O42").Value = "10.9 °C"
$ws.Range("E43").Value = "2026-02-21 20:20:21"
$ws.Range("O43").Value = "7.2 °C"
$ws.Range("E44").Value = "2026-02-21 20:20:23"
$ws.Range("N44").Value = "-0.6 °C 19:46 TU"
$ws.Range("E45").Value = "2026-02-21 20:20:26"
$ws.Range("E46").Value = "2026-02-21 20:20:29"
$ws.Range("O46").Value = "10.0 °C"
